# Applies the "Updated cryptos list" price/volume refresh to the worksheet.
# Column D = Price, Column E = Volume(1h); both stored as plain text strings
# (inline strings in the source workbook), matching the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.903.59'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '2.821.62'
$ws.Range('E3').Value = '  +1.65%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''354.65'
$ws.Range('E5').Value = '  +6.28%  '
$ws.Range('D6').Value = '''113.48'
$ws.Range('E6').Value = '  -2.81%  '
$ws.Range('E7').Value = '  +4.42%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '''0.608'
$ws.Range('E9').Value = '  +5.47%  '
$ws.Range('E10').Value = '  -1.43%  '
$ws.Range('D11').Value = '''0.0853'
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  +1.11%  '
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').Value = '''7.77'
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '3.268.61'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '2.819.38'
$ws.Range('E16').Value = '  +1.09%  '
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '51.834.25'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').Value = '''7.49'
$ws.Range('E19').Value = '  +8.96%  '
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('D21').Value = '''13.47'
$ws.Range('E21').Value = '  -0.65%  '
$ws.Range('D22').Value = '0.0₃0988'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('D23').Value = '''270.64'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = '''69.76'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').Value = '''26.79'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('D27').Value = '''0.999'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '''10.32'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('E29').Value = '  +0.83%  '
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').Value = '''50.81'
$ws.Range('E31').Value = '  +1.07%  '
$ws.Range('D32').Value = '''34.05'
$ws.Range('E32').Value = '  -3.69%  '
$ws.Range('D33').Value = '''0.0455'
$ws.Range('E33').Value = '  +26.91%  '
$ws.Range('D34').Value = '''5.86'
$ws.Range('D35').Value = '''5.32'
$ws.Range('E35').Value = '  +6.38%  '
$ws.Range('D36').Value = '''0.0828'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('D39').Value = '''3.23'
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  -4.92%  '
$ws.Range('D41').Value = '''23.91'
$ws.Range('E41').Value = '  +2.12%  '
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').Value = '''125.72'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').Value = '''2.29'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('D46').Value = '''3.39'
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('D47').Value = '2.095.00'
$ws.Range('E47').Value = '  +0.28%  '
$ws.Range('D48').Value = '''2.27'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('E50').Value = '  +5.60%  '
$ws.Range('D51').Value = '''60.81'
$ws.Range('E51').Value = '  +0.02%  '
